# Update crypto price/volume table with latest scraped values (GitHub Actions refresh).
# Rows 13/14 also swap their coin (Chainlink <-> WrappedEther) to reflect new ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.525.32'
$ws.Range('E2').Value = '  -0.23%  '

$ws.Range('D3').Value = '1.810.53'
$ws.Range('E3').Value = '  +0.63%  '

$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.47%  '

$ws.Range('E6').Value = '  +4.04%  '

$ws.Range('E7').Value = '  +0.08%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '34.94'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.35%  '

$ws.Range('E9').Value = '  +2.52%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0698'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.34%  '

$ws.Range('E11').Value = '  +0.79%  '

$ws.Range('D12').Value = '2.072.29'
$ws.Range('E12').Value = '  +0.57%  '

$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.21'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.69%  '

$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.816.00'
$ws.Range('E14').Value = '  +1.13%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.650'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.54%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.53'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.84%  '

$ws.Range('D17').Value = '34.516.18'
$ws.Range('E17').Value = '  -0.25%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.29'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.57%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '246.54'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.21%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.45'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.46%  '

$ws.Range('E22').Value = '  +0.10%  '

$ws.Range('E23').Value = '  +0.15%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '172.91'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.18%  '

$ws.Range('E25').Value = '  +2.20%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.14'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +10.99%  '

$ws.Range('E27').Value = '  +1.09%  '

$ws.Range('E28').Value = '  +2.55%  '

$ws.Range('E29').Value = '  -0.03%  '

$ws.Range('E30').Value = '  -0.21%  '

$ws.Range('E31').Value = '  +2.06%  '

$ws.Range('E32').Value = '  +1.85%  '

$ws.Range('E33').Value = '  +0.42%  '

$ws.Range('E34').Value = '  +0.60%  '

$ws.Range('D35').Value = '1.397.91'

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.681'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.56%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.49'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.26%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.08'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.44%  '

$ws.Range('E39').Value = '  -0.13%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '83.91'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.09%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.965'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.75%  '

$ws.Range('E42').Value = '  +2.84%  '

$ws.Range('E43').Value = '  -0.17%  '

$ws.Range('E44').Value = '  +5.71%  '

$ws.Range('E45').Value = '  -3.75%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0513'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.27%  '

$ws.Range('E47').Value = '  -1.46%  '

$ws.Range('D48').Value = '1.971.92'
$ws.Range('E48').Value = '  +0.56%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '105.39'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.15%  '

$ws.Range('E50').Value = '  +2.13%  '

$ws.Range('E51').Value = '  +0.11%  '
